$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix cell TYPE changes first (number <-> text) by copying a template cell
# that already has the desired style+type, then (if needed) overwrite the value. ---

# number -> text "0" (style 14, shared string "0")
$ws.Range("G30").Copy($ws.Range("D15"))
$ws.Range("G30").Copy($ws.Range("C26"))
$ws.Range("G30").Copy($ws.Range("D26"))
$ws.Range("G30").Copy($ws.Range("C27"))
$ws.Range("G30").Copy($ws.Range("D27"))
$ws.Range("G30").Copy($ws.Range("F28"))
$ws.Range("G30").Copy($ws.Range("F29"))
$ws.Range("G30").Copy($ws.Range("F30"))

# number -> text "***.*" (style 14, shared string "***.*")
$ws.Range("H30").Copy($ws.Range("E15"))
$ws.Range("H30").Copy($ws.Range("E26"))
$ws.Range("H30").Copy($ws.Range("E27"))

# text -> number (style 16, plain integer format)
$ws.Range("F15").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 4

# text -> number (style 15, percentage-decimal format)
$ws.Range("N14").Copy($ws.Range("E18"))
$ws.Range("E18").Value = 0

# --- Plain value updates (style/type unchanged) ---
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -80
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 133.333333333333
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 18
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = 38.461538461538
$ws.Range("L16").Value = 63.636363636363
$ws.Range("M16").Value = -48.571428571428
$ws.Range("N16").Value = -82.692307692307
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 10.526315789473
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 24
$ws.Range("K17").Value = 33.333333333333
$ws.Range("L17").Value = 68.421052631578
$ws.Range("M17").Value = 146.153846153846
$ws.Range("N17").Value = -5.882352941176
$ws.Range("C18").Value = 4
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = 45.454545454545
$ws.Range("L18").Value = 45.454545454545
$ws.Range("M18").Value = -30.434782608695
$ws.Range("N18").Value = -88.652482269503
$ws.Range("C19").Value = 14
$ws.Range("E19").Value = -6.666666666666
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = -17.647058823529
$ws.Range("I19").Value = 66
$ws.Range("J19").Value = 81
$ws.Range("K19").Value = -18.518518518518
$ws.Range("L19").Value = 57.142857142857
$ws.Range("M19").Value = 120
$ws.Range("N19").Value = 13.793103448275
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 4
$ws.Range("I20").Value = 32
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = 10.344827586206
$ws.Range("L20").Value = 128.571428571429
$ws.Range("M20").Value = -5.882352941176
$ws.Range("N20").Value = -90.934844192634
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 131
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = -2.238805970149
$ws.Range("I21").Value = 165
$ws.Range("J21").Value = 160
$ws.Range("K21").Value = 3.125
$ws.Range("L21").Value = 70.103092783505
$ws.Range("M21").Value = 20.437956204379
$ws.Range("N21").Value = -76.361031518624
$ws.Range("C22").Value = 3
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = 33.333333333333
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = 100
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = 19.753086419753
$ws.Range("I24").Value = 111
$ws.Range("J24").Value = 96
$ws.Range("K24").Value = 15.625
$ws.Range("L24").Value = 50
$ws.Range("M24").Value = 85
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 37.037037037037
$ws.Range("I25").Value = 47
$ws.Range("J25").Value = 35
$ws.Range("K25").Value = 34.285714285714
$ws.Range("L25").Value = 51.612903225806
$ws.Range("M25").Value = 11.904761904761
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = 60
$ws.Range("L27").Value = 42.857142857142
$ws.Range("H28").Value = -100
$ws.Range("N28").Value = -60
$ws.Range("H29").Value = -100
$ws.Range("N29").Value = -80

# --- Header text edits (edit specific characters within the rich-text run,
# preserving the surrounding run formatting/text) ---
$ws.Range("A8").Characters(21, 1).Text = "5"
$ws.Range("C9").Characters(27, 9).Text = "1/30/2023"
$ws.Range("C9").Characters(47, 9).Text = "2/5/2023"
